$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 17.571428
$ws.Range("I11").Value = 17.571428
$ws.Range("K11").Value = 17.571428
$ws.Range("M11").Value = 122.428572
$ws.Range("H74").Value = 2930.0889
$ws.Range("I74").Value = 3053.7334
$ws.Range("J74").Value = 2868.2666
$ws.Range("K74").Value = 3053.7334
$ws.Range("L74").Value = 2868.2666
$ws.Range("M74").Value = -2117.7334
$ws.Range("N74").Value = -4740.2666
$ws.Range("H76").Value = 3370.9443
$ws.Range("I76").Value = 2979.9092
$ws.Range("J76").Value = 3985.4285
$ws.Range("K76").Value = 2979.9092
$ws.Range("L76").Value = 3985.4285
$ws.Range("M76").Value = -2664.9092
$ws.Range("N76").Value = -4615.4285
$ws.Range("H77").Value = 2930.0889
$ws.Range("I77").Value = 3053.7334
$ws.Range("J77").Value = 2868.2666
$ws.Range("K77").Value = 15268.667
$ws.Range("L77").Value = 14341.333
$ws.Range("M77").Value = -10588.667
$ws.Range("N77").Value = -23701.333
$ws.Range("H79").Value = 3370.9443
$ws.Range("I79").Value = 2979.9092
$ws.Range("J79").Value = 3985.4285
$ws.Range("K79").Value = 2979.9092
$ws.Range("L79").Value = 3985.4285
$ws.Range("M79").Value = -1887.9092
$ws.Range("N79").Value = -6169.4285
$ws.Range("H112").Value = 1364.5
$ws.Range("I112").Value = 812.4
$ws.Range("J112").Value = 1576.8462
$ws.Range("K112").Value = 2437.2
$ws.Range("L112").Value = 4730.5386
$ws.Range("M112").Value = -1329.2
$ws.Range("N112").Value = -6946.5386
$ws.Range("H116").Value = 1896058.9
$ws.Range("I116").Value = 2254298
$ws.Range("J116").Value = 2508.8572
$ws.Range("K116").Value = 2254298
$ws.Range("L116").Value = 2508.8572
$ws.Range("M116").Value = -2250856
$ws.Range("N116").Value = -9392.8572
$ws.Range("H127").Value = 494.63635
$ws.Range("I127").Value = 494.63635
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 1483.90905
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 3476.09095
$ws.Range("N127").Value = ""
$ws.Range("H129").Value = 910.931
$ws.Range("I129").Value = 359.625
$ws.Range("J129").Value = 1120.9524
$ws.Range("K129").Value = 1078.875
$ws.Range("L129").Value = 3362.857199999999
$ws.Range("M129").Value = 3921.125
$ws.Range("N129").Value = -13362.8572
$ws.Range("H137").Value = 1331.3864
$ws.Range("I137").Value = 1224.3077
$ws.Range("J137").Value = 1486.0555
$ws.Range("K137").Value = 3672.9231
$ws.Range("L137").Value = 4458.166499999999
$ws.Range("M137").Value = -1122.9231
$ws.Range("N137").Value = -9558.166499999999
$ws.Range("H138").Value = 2656.4849
$ws.Range("I138").Value = 1399.25
$ws.Range("J138").Value = 3058.8
$ws.Range("K138").Value = 4197.75
$ws.Range("L138").Value = 9176.400000000001
$ws.Range("M138").Value = 942.25
$ws.Range("N138").Value = -19456.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2132.8823
$ws.Range("I45").Value = 2404.5
$ws.Range("J45").Value = 1744.8572
$ws.Range("K45").Value = 2404.5
$ws.Range("L45").Value = 1744.8572
$ws.Range("M45").Value = -2027.5
$ws.Range("N45").Value = -2498.8572
$ws.Range("H61").Value = 1516.775
$ws.Range("I61").Value = 1215.1936
$ws.Range("J61").Value = 2555.5557
$ws.Range("K61").Value = 1215.1936
$ws.Range("L61").Value = 2555.5557
$ws.Range("M61").Value = -1003.1936
$ws.Range("N61").Value = -2979.5557
$ws.Range("H74").Value = 2724.0857
$ws.Range("I74").Value = 3063.842
$ws.Range("J74").Value = 2320.625
$ws.Range("K74").Value = 3063.842
$ws.Range("L74").Value = 2320.625
$ws.Range("M74").Value = -2189.842
$ws.Range("N74").Value = -4068.625
$ws.Range("H77").Value = 2724.0857
$ws.Range("I77").Value = 3063.842
$ws.Range("J77").Value = 2320.625
$ws.Range("K77").Value = 15319.21
$ws.Range("L77").Value = 11603.125
$ws.Range("M77").Value = -10951.21
$ws.Range("N77").Value = -20339.125
$ws.Range("H121").Value = 45255
$ws.Range("J121").Value = 45255
$ws.Range("L121").Value = 45255
$ws.Range("N121").Value = -48749
$ws.Range("H136").Value = 1516.775
$ws.Range("I136").Value = 1215.1936
$ws.Range("J136").Value = 2555.5557
$ws.Range("K136").Value = 3645.5808
$ws.Range("L136").Value = 7666.6671
$ws.Range("M136").Value = -1095.5808
$ws.Range("N136").Value = -12766.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5908.3335
$ws.Range("I86").Value = 6090
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 6090
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -4967
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 5908.3335
$ws.Range("I89").Value = 6090
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 30450
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -24834
$ws.Range("N89").Value = -36232
$ws.Range("H99").Value = 3860.5417
$ws.Range("I99").Value = 3793
$ws.Range("J99").Value = 4333.3335
$ws.Range("K99").Value = 3793
$ws.Range("L99").Value = 4333.3335
$ws.Range("M99").Value = -2295
$ws.Range("N99").Value = -7329.3335
$ws.Range("H126").Value = 3860.5417
$ws.Range("I126").Value = 3793
$ws.Range("J126").Value = 4333.3335
$ws.Range("K126").Value = 11379
$ws.Range("L126").Value = 13000.0005
$ws.Range("M126").Value = -8909
$ws.Range("N126").Value = -17940.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1013.4545
$ws.Range("J131").Value = 1052.525
$ws.Range("L131").Value = 3157.575
$ws.Range("N131").Value = -13237.575

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 216.8
$ws.Range("I2").Value = 139.6
$ws.Range("J2").Value = 294
$ws.Range("K2").Value = 139.6
$ws.Range("L2").Value = 294
$ws.Range("M2").Value = -26.59999999999999
$ws.Range("N2").Value = -520

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7114.6523
$ws.Range("I136").Value = 17212.75
$ws.Range("J136").Value = 1729
$ws.Range("K136").Value = 51638.25
$ws.Range("L136").Value = 5187
$ws.Range("M136").Value = -49088.25
$ws.Range("N136").Value = -10287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 198.96
$ws.Range("I113").Value = 193.05
$ws.Range("J113").Value = 222.6
$ws.Range("K113").Value = 579.1500000000001
$ws.Range("L113").Value = 667.8
$ws.Range("M113").Value = 1590.85
$ws.Range("N113").Value = -5007.8
$ws.Range("H122").Value = 1333.5652
$ws.Range("I122").Value = 1074.75
$ws.Range("J122").Value = 1471.6
$ws.Range("K122").Value = 3224.25
$ws.Range("L122").Value = 4414.799999999999
$ws.Range("M122").Value = -774.25
$ws.Range("N122").Value = -9314.799999999999
$ws.Range("H133").Value = 38845
$ws.Range("J133").Value = 38845
$ws.Range("L133").Value = 38845
$ws.Range("N133").Value = -48965
$ws.Range("H136").Value = 1371.4242
$ws.Range("I136").Value = 724.04346
$ws.Range("J136").Value = 2860.4
$ws.Range("K136").Value = 2172.13038
$ws.Range("L136").Value = 8581.200000000001
$ws.Range("M136").Value = 377.8696199999999
$ws.Range("N136").Value = -13681.2
